$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.077.66'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +0.73%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.891.79'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +1.52%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.57'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +0.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  +0.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5183'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  +2.47%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3761'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  +3.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07220'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +0.67%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.19'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  +2.51%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9022'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  +1.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07649'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  +1.84%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.892.51'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  +1.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.53'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  -0.25%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.243'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +0.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +0.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008513'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +0.02%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.45'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  +1.73%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9998'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  +0.04%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.126.38'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +0.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.057'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +0.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.078.58'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  -0.94%  '

$ws.Range("E23").Value = '  +2.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.385'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  -0.35%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.320'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  +11.67%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '145.58'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  -1.61%  '

$ws.Range("B27").Value = 'EthereumClassic'

$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.09'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  +1.17%  '

$ws.Range("B28").Value = 'Toncoin'

$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.724'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  -3.10%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.39'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  +1.07%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.922'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  +5.50%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.799'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  +2.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09202'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  +0.16%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05054'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  -1.52%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.248'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +8.17%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7715'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +2.59%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.985'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  +0.55%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.283'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  +1.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.591'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +0.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5614'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +0.74%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01993'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  -0.23%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.071'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +0.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.058'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +5.79%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.630'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  +0.70%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '118.57'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +2.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1510'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +2.67%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4841'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  +3.35%  '

$ws.Range("B47").Value = 'EnergySwap'

$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.16'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  +1.05%  '

$ws.Range("B48").Value = 'PaxDollar'

$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9998'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  +0.17%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.599'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +2.71%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.74'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +2.94%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.01'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +1.66%  '
